# Updated cryptos list on Tue Sep  3 02:48:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "59.454.42"
$ws.Cells.Item(2, 5).Value = "  +3.72%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.533.64"
$ws.Cells.Item(3, 5).Value = "  +4.50%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "534.78"
$ws.Cells.Item(5, 5).Value = "  +4.35%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "136.40"
$ws.Cells.Item(6, 5).Value = "  +5.72%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.40%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.569"
$ws.Cells.Item(8, 5).Value = "  +3.85%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.531.02"
$ws.Cells.Item(9, 5).Value = "  +3.92%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +5.18%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.88%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "5.25"
$ws.Cells.Item(12, 5).Value = "  +1.54%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.337"
$ws.Cells.Item(13, 5).Value = "  +1.55%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.979.92"
$ws.Cells.Item(14, 5).Value = "  +4.36%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "59.399.24"
$ws.Cells.Item(15, 5).Value = "  +3.83%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "22.59"
$ws.Cells.Item(16, 5).Value = "  +4.14%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +4.59%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.533.74"
$ws.Cells.Item(18, 5).Value = "  +4.21%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.80"
$ws.Cells.Item(19, 5).Value = "  +3.65%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.26"
$ws.Cells.Item(20, 5).Value = "  +4.08%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "324.53"
$ws.Cells.Item(21, 5).Value = "  +3.19%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.15"
$ws.Cells.Item(22, 5).Value = "  +8.96%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.999"
$ws.Cells.Item(23, 5).Value = "  +0.14%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "65.85"
$ws.Cells.Item(24, 5).Value = "  +3.48%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.411"
$ws.Cells.Item(25, 5).Value = "  +1.35%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.999"
$ws.Cells.Item(26, 5).Value = "  +0.15%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "Kaspa"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.162"
$ws.Cells.Item(27, 5).Value = "  +2.24%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +5.93%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +7.32%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "174.40"
$ws.Cells.Item(30, 5).Value = "  +3.44%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.75"
$ws.Cells.Item(31, 5).Value = "  +5.75%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Fetch.AI"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.23"
$ws.Cells.Item(32, 5).Value = "  +5.84%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.43"
$ws.Cells.Item(33, 5).Value = "  +3.24%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.999"
$ws.Cells.Item(34, 5).Value = "  +0.02%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.00"
$ws.Cells.Item(35, 5).Value = "  +0.29%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "18.33"
$ws.Cells.Item(36, 5).Value = "  +3.69%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.27"
$ws.Cells.Item(37, 5).Value = "  +0.21%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.04"
$ws.Cells.Item(38, 5).Value = "  +3.89%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.53"
$ws.Cells.Item(39, 5).Value = "  +5.97%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "36.79"
$ws.Cells.Item(40, 5).Value = "  +1.93%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.794"
$ws.Cells.Item(41, 5).Value = "  +2.62%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "283.51"
$ws.Cells.Item(42, 5).Value = "  +5.97%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Filecoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.53"
$ws.Cells.Item(43, 5).Value = "  +5.06%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.17"
$ws.Cells.Item(44, 5).Value = "  +5.96%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "133.17"
$ws.Cells.Item(45, 5).Value = "  +10.92%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.602"
$ws.Cells.Item(46, 5).Value = "  +3.02%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0932"
$ws.Cells.Item(47, 5).Value = "  +3.04%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0513"
$ws.Cells.Item(48, 5).Value = "  +6.28%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0221"
$ws.Cells.Item(49, 5).Value = "  +5.74%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "17.30"
$ws.Cells.Item(50, 5).Value = "  +5.41%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "1.768.38"
$ws.Cells.Item(51, 5).Value = "  +4.17%  "
